# Fruta / hortaliza, semanal
# Rows 3-10 are being re-shuffled: the data that previously lived in one
# row (columns D, J, K, L, M, P) now lives in a different row, per the
# permutation observed in the diff. Columns A, B, C, E, F, G, H, I, N, O, Q, R
# are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per destination row, taken from the corresponding source row
# in the original data (D, J, K, L, M, P)
$values = @{
    3  = @{ D = 44977; J = 400; K = 16500; L = 17000; M = 16750; P = 931 }
    4  = @{ D = 44557; J = 400; K = 13000; L = 14000; M = 13500; P = 750 }
    5  = @{ D = 44984; J = 200; K = 17000; L = 18000; M = 17500; P = 972 }
    6  = @{ D = 44547; J = 200; K = 13000; L = 14000; M = 13500; P = 750 }
    7  = @{ D = 44998; J = 320; K = 17000; L = 18000; M = 17500; P = 972 }
    8  = @{ D = 44957; J = 400; K = 21000; L = 22000; M = 21500; P = 1194 }
    9  = @{ D = 44568; J = 500; K = 15000; L = 16000; M = 15500; P = 861 }
    10 = @{ D = 44960; J = 400; K = 19500; L = 20000; M = 19750; P = 1097 }
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Range("D$row").Value = $v.D
    $ws.Range("J$row").Value = $v.J
    $ws.Range("K$row").Value = $v.K
    $ws.Range("L$row").Value = $v.L
    $ws.Range("M$row").Value = $v.M
    $ws.Range("P$row").Value = $v.P
}
